$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.896.37"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.640.83"
$ws.Range("E3").Value = "  -1.49%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.28%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "213.45"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.46%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.5205"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.20%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2600"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.14%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06321"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -0.07%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "20.57"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -1.94%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07673"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").Value = "1.630.34"
$ws.Range("E12").Value = "  -2.26%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "4.410"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("D14").Value = "1.863.54"
$ws.Range("E14").Value = "  -1.59%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.5480"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +2.01%  "
$ws.Range("D16").Value = "0.0₅8190"
$ws.Range("E16").Value = "  +2.87%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "64.39"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("D18").Value = "25.898.50"
$ws.Range("E18").Value = "  -0.99%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "4.682"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.52%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "188.41"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.50%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.13"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.89%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "6.239"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.57%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.18%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "143.09"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -3.98%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.1236"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.87%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "7.347"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.85%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "15.84"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.98%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "1.409"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.52%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.05903"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -4.51%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.256"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.78%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.386"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.37%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "3.391"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.43%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.635"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.01%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.9863"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.40%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.397"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.24%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.735"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.77%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.5590"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -5.16%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.01596"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.11%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "5.832"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -3.00%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.8495"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.08%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").Value = "1.018.54"
$ws.Range("E43").Value = "  -8.04%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "98.54"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.37%  "
$ws.Range("D45").Value = "1.788.32"
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("D46").Value = "0.0₈111"
$ws.Range("E46").Value = "  +0.38%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "55.37"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.07%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.14%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "8.026"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.55%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.05138"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -2.11%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.4208"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.15%  "
